# Update "想去人数" (want-to-go count) values in columns F on the
# "展览" and "全部类型" worksheets to match the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1640
$ws1.Range("F6").Value = 3247
$ws1.Range("F7").Value = 794
$ws1.Range("F8").Value = 2020
$ws1.Range("F9").Value = 1937
$ws1.Range("F11").Value = 347
$ws1.Range("F13").Value = 1597
$ws1.Range("F17").Value = 38
$ws1.Range("F18").Value = 1414
$ws1.Range("F19").Value = 508
$ws1.Range("F20").Value = 610
$ws1.Range("F22").Value = 10593
$ws1.Range("F23").Value = 9757
$ws1.Range("F24").Value = 840
$ws1.Range("F26").Value = 1822
$ws1.Range("F27").Value = 138
$ws1.Range("F28").Value = 401

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 1640
$ws4.Range("F8").Value = 3247
$ws4.Range("F9").Value = 794
$ws4.Range("F10").Value = 2020
$ws4.Range("F11").Value = 1937
$ws4.Range("F13").Value = 347
$ws4.Range("F15").Value = 1597
$ws4.Range("F21").Value = 38
$ws4.Range("F22").Value = 1414
$ws4.Range("F23").Value = 508
$ws4.Range("F24").Value = 610
$ws4.Range("F26").Value = 10593
$ws4.Range("F27").Value = 9757
$ws4.Range("F28").Value = 840
$ws4.Range("F30").Value = 1822
$ws4.Range("F33").Value = 138
$ws4.Range("F34").Value = 401
